$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 256, pushing existing rows 256-266 down to 257-267
$ws.Rows.Item(256).Insert()

# Fill in the new row 256 with the weekly price record
$ws.Range("A256").Value = 11
$ws.Range("B256").Value = "Vega Monumental Concepción"
$ws.Range("C256").Value = "Bíobío"
$ws.Range("D256").Value = 45041
$ws.Range("E256").Value = 8
$ws.Range("F256").Value = "Fruta"
$ws.Range("G256").Value = 100108
$ws.Range("H256").Value = "Tropicales y subtropicales"
$ws.Range("I256").Value = 100108005
$ws.Range("J256").Value = "Piña"
$ws.Range("K256").Value = "Caramelo"
$ws.Range("L256").Value = "Segunda"
$ws.Range("M256").Value = 100
$ws.Range("N256").Value = 16000
$ws.Range("O256").Value = 17000
$ws.Range("P256").Value = 16500
$ws.Range("Q256").Value = "$/caja 14 unidades"
$ws.Range("R256").Value = "Ecuador"
$ws.Range("S256").Value = 1179
$ws.Range("T256").Value = 14
